# Weekly refresh: a new daily price record (2023-09-08) was added to the
# top of the data block, pushing every existing record down by one row.
# Net effect on the sheet: insert one row at row 29 (the first data row
# for this market/category block), fill it with the new record, and the
# used range grows from A1:R116 to A1:R117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (rows 29:116) down by inserting a new row at 29.
# Excel's Insert copies formatting down from the row above, so column D
# keeps its date style automatically.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(29, 3).Value = "Metropolitana"
$ws.Cells.Item(29, 4).Value = 45177
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = 100112035
$ws.Cells.Item(29, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 410
$ws.Cells.Item(29, 11).Value = 17000
$ws.Cells.Item(29, 12).Value = 18000
$ws.Cells.Item(29, 13).Value = 17439
$ws.Cells.Item(29, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(29, 16).Value = 1163
$ws.Cells.Item(29, 17).Value = 15
$ws.Cells.Item(29, 18).Value = "Hortaliza"
